$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("2021-07-27", "overview", "K02000001", "United Kingdom", 5745526, 23511, 131, 129303)
    ,@("2021-07-28", "overview", "K02000001", "United Kingdom", 5770928, 27734, 91, 129430)
    ,@("2021-07-29", "overview", "K02000001", "United Kingdom", 5801561, 31117, 85, 129515)
    ,@("2021-07-30", "overview", "K02000001", "United Kingdom", 5830774, 29622, 68, 129583)
    ,@("2021-07-31", "overview", "K02000001", "United Kingdom", 5856528, 26144, 71, 129654)
    ,@("2021-08-01", "overview", "K02000001", "United Kingdom", 5880667, 24470, 65, 129719)
    ,@("2021-08-02", "overview", "K02000001", "United Kingdom", 5902354, 21952, 24, 129743)
    ,@("2021-08-03", "overview", "K02000001", "United Kingdom", 5923820, 21691, 138, 129881)
    ,@("2021-08-04", "overview", "K02000001", "United Kingdom", 5952756, 29312, 119, 130000)
    ,@("2021-08-05", "overview", "K02000001", "United Kingdom", 5982581, 30215, 86, 130086)
    ,@("2021-08-06", "overview", "K02000001", "United Kingdom", 6014023, 31808, 92, 130178)
    ,@("2021-08-07", "overview", "K02000001", "United Kingdom", 6042252, 28612, 103, 130281)
    ,@("2021-08-08", "overview", "K02000001", "United Kingdom", 6069362, 27429, 39, 130320)
    ,@("2021-08-09", "overview", "K02000001", "United Kingdom", 6094243, 25161, 37, 130357)
    ,@("2021-08-10", "overview", "K02000001", "United Kingdom", 6117540, 23510, 146, 130503)
    ,@("2021-08-11", "overview", "K02000001", "United Kingdom", 6146800, 29612, 104, 130607)
    ,@("2021-08-12", "overview", "K02000001", "United Kingdom", 6179506, 33074, 94, 130701)
    ,@("2021-08-13", "overview", "K02000001", "United Kingdom", 6211868, 32700, 100, 130801)
)

$startRow = 350
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    # Leading apostrophe forces text interpretation so the date-shaped
    # string (e.g. 2021-07-27) is not auto-converted to a date serial,
    # matching the source data's plain-text date column.
    $ws.Cells.Item($r, 1).Value = "'" + $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}
